# Auto-generated Excel COM-interop edit script
# Applies updates to columns E, F, G (selected rows) and adds column I ("Other found locations")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-B64Cell($addr, $b64) {
    $bytes = [System.Convert]::FromBase64String($b64)
    $text = [System.Text.Encoding]::UTF8.GetString($bytes)
    $ws.Range($addr).Value = $text
}

Set-B64Cell "I1" "T3RoZXIgZm91bmQgbG9jYXRpb25z"
Set-B64Cell "E2" "W0FrZGlzJUNlem1pIEElY29yZUdpdmVzTm9FbWFpbCUwLCAgIENhbyVZaS15dWFuJWNvcmVHaXZlc05vRW1haWwlMCwgICBEb25nJVhpYW5nJWNvcmVHaXZlc05vRW1haWwlMCwgICBHYW8lWWEtZG9uZyVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgWWFuJVlvdS1xaW4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIFlhbmclWWktYmluJWNvcmVHaXZlc05vRW1haWwlMCwgICBZdWFuJVlhLWRvbmclY29yZUdpdmVzTm9FbWFpbCUwLCAgIFpoYW5nJUppbi1qaW4lY29yZUdpdmVzTm9FbWFpbCUwXQ=="
Set-B64Cell "F2" "bm90IGZvdW5k"
Set-B64Cell "G2" "Ti9B"
Set-B64Cell "I2" ""
Set-B64Cell "E3" "W1pob25nbGlhbmclV2FuZyVOVUxMJTAsICAgICBCb2hhbiVZYW5nJU5VTEwlMCwgICAgIFFpYW53ZW4lTGklTlVMTCUwLCAgICAgTHUlV2VuJU5VTEwlMCwgICAgIFJ1aWd1YW5nJVpoYW5nJXpyZzI3QDE2My5jb20lMF0="
Set-B64Cell "I3" "X1BNQw=="
Set-B64Cell "E4" "W1hpYW9ibyVZYW5nJU5VTEwlMCwgICAgIFl1YW4lWXUlTlVMTCUwLCAgICAgSmlxaWFuJVh1JU5VTEwlMCwgICAgIEh1YXFpbmclU2h1JU5VTEwlMCwgICAgIEppYSdhbiVYaWElTlVMTCUwLCAgICAgSG9uZyVMaXUlTlVMTCUwLCAgICAgWW9uZ3JhbiVXdSVOVUxMJTAsICAgICBMdSVaaGFuZyVOVUxMJTAsICAgICBaaHVpJVl1JU5VTEwlMCwgICAgIE1pbmdoYW8lRmFuZyVOVUxMJTAsICAgICBUaW5nJVl1JU5VTEwlMCwgICAgIFlheGluJVdhbmclTlVMTCUwLCAgICAgU2hhbmd3ZW4lUGFuJU5VTEwlMCwgICAgIFhpYW9qaW5nJVpvdSVOVUxMJTAsICAgICBTaGl5aW5nJVl1YW4lTlVMTCUwLCAgICAgWW91JVNoYW5nJU5VTEwlMF0="
Set-B64Cell "I4" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "I5" ""
Set-B64Cell "E6" "W1lpaHVpJUh1YW5nJU5VTEwlMCwgICAgIE1lbmdxaSVUdSVOVUxMJTAsICAgICBTaGlwZWklV2FuZyVOVUxMJTAsICAgICBTaWNoYW8lQ2hlbiVOVUxMJTAsICAgICBXZWklWmhvdSVOVUxMJTAsICAgICBEYW55YW5nJUNoZW4lTlVMTCUwLCAgICAgTGluJVpob3UlTlVMTCUwLCAgICAgTWluJVdhbmclTlVMTCUwLCAgICAgWWFuJVpoYW8lTlVMTCUwLCAgICAgV2VuJVplbmclTlVMTCUwLCAgICAgUWklSHVhbmclTlVMTCUwLCAgICAgSGFpJ2JvJVh1JU5VTEwlMCwgICAgIFplbWluZyVMaXUlTlVMTCUwLCAgICAgTGlhbmclR3VvJU5VTEwlMF0="
Set-B64Cell "I6" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E7" "W0t1bmh1YSVMaSVOVUxMJTAsICAgICBKaW9uZyVXdSVOVUxMJTAsICAgICBGYXFpJVd1JU5VTEwlMCwgICAgIERhamluZyVHdW8lTlVMTCUwLCAgICAgTGlubGklQ2hlbiVOVUxMJTAsICAgICBaaGVuZyVGYW5nJU5VTEwlMCwgICAgIENodWFubWluZyVMaSVOVUxMJTBd"
Set-B64Cell "I7" "X1BNQw=="
Set-B64Cell "E8" "W1hpJVh1JU5VTEwlMCwgICAgIENoZW5nY2hlbmclWXUlTlVMTCUwLCAgICAgSmluZyVRdSVOVUxMJTAsICAgICBMaWVndWFuZyVaaGFuZyVOVUxMJTAsICAgICBTb25nZmVuZyVKaWFuZyVOVUxMJTAsICAgICBEZXlhbmclSHVhbmclTlVMTCUwLCAgICAgQmlodWElQ2hlbiVOVUxMJTAsICAgICBaaGlwaW5nJVpoYW5nJU5VTEwlMCwgICAgIFdhbmh1YSVHdWFuJU5VTEwlMCwgICAgIFpob3VrdW4lTGluZyVOVUxMJTAsICAgICBSdWklSmlhbmclTlVMTCUwLCAgICAgVGlhbmxpJUh1JU5VTEwlMCwgICAgIFlhbiVEaW5nJU5VTEwlMCwgICAgIExpbiVMaW4lTlVMTCUwLCAgICAgUWluZ3hpbiVHYW4lTlVMTCUwLCAgICAgTGlhbmdwaW5nJUx1byV0bHVvbHBAam51LmVkdS5jbiUwLCAgICAgWGlhb3BpbmclVGFuZyV4dGFuZ0AyMWNuLmNvbSUwLCAgICAgSmlueGluJUxpdSVMaXVqeDgzNzEwMzc4QDEyNi5jb20lMF0="
Set-B64Cell "I8" "X1BNQ19TcHJpbmdlcg=="
Set-B64Cell "E9" "W1dlbmppZSVZYW5nJU5VTEwlMCwgICAgIFFpcWklQ2FvJU5VTEwlMCwgICAgIExlJVFpbiVOVUxMJTAsICAgICBYaWFveWFuZyVXYW5nJU5VTEwlMCwgICAgIFplbmdodWklQ2hlbmclTlVMTCUwLCAgICAgQXNoYW4lUGFuJU5VTEwlMCwgICAgIEppYW55aSVEYWklTlVMTCUwLCAgICAgUWluZ2ZlbmclU3VuJU5VTEwlMCwgICAgIEZlbmdxdWFuJVpoYW8lTlVMTCUwLCAgICAgSmllbWluZyVRdSVOVUxMJTAsICAgICBGdWh1YSVZYW4lTlVMTCUwXQ=="
Set-B64Cell "I9" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E10" "W1dlbmppZSVZYW5nJU5VTEwlMCwgICAgIFFpcWklQ2FvJU5VTEwlMCwgICAgIExlJVFpbiVOVUxMJTAsICAgICBYaWFveWFuZyVXYW5nJU5VTEwlMCwgICAgIFplbmdodWklQ2hlbmclTlVMTCUwLCAgICAgQXNoYW4lUGFuJU5VTEwlMCwgICAgIEppYW55aSVEYWklTlVMTCUwLCAgICAgUWluZ2ZlbmclU3VuJU5VTEwlMCwgICAgIEZlbmdxdWFuJVpoYW8lTlVMTCUwLCAgICAgSmllbWluZyVRdSVOVUxMJTAsICAgICBGdWh1YSVZYW4lTlVMTCUwXQ=="
Set-B64Cell "I10" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E11" "W1NoYW9xaW5nJUxlaSVOVUxMJTAsICAgICBGYW5nJUppYW5nJU5VTEwlMCwgICAgIFdhdGluZyVTdSVOVUxMJTAsICAgICBDaGFuZyVDaGVuJU5VTEwlMCwgICAgIEppbmdsaSVDaGVuJU5VTEwlMCwgICAgIFdlaSVNZWklTlVMTCUwLCAgICAgTGktWWluZyVaaGFuJU5VTEwlMCwgICAgIFlpZmFuJUppYSVOVUxMJTAsICAgICBMaWFuZ3FpbmclWmhhbmclTlVMTCUwLCAgICAgRGFueW9uZyVMaXUlTlVMTCUwLCAgICAgWmhvbmctWXVhbiVYaWEleGlhemhvbmd5dWFuMjAwNUBhbGl5dW4uY29tJTAsICAgICBaaGVuZ3l1YW4lWGlhJXp5eGlhQGhrdWNjLmhrdS5oayUwXQ=="
Set-B64Cell "I11" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E12" "W1l1biVGZW5nJU5VTEwlMCwgICAgIFl1biVMaW5nJU5VTEwlMCwgICAgIFl1biVMaW5nJU5VTEwlMCwgICAgIFRhbyVCYWklTlVMTCUzLCAgICAgVGFvJUJhaSVOVUxMJTAsICAgICBZdXNhbmclWGllJU5VTEwlMSwgICAgIEppZSVIdWFuZyVOVUxMJTIsICAgICBKaWUlSHVhbmclTlVMTCUwLCAgICAgSmlhbiVMaSVOVUxMJTEsICAgICBXZWluaW5nJVhpb25nJU5VTEwlMSwgICAgIERleGlhbmclWWFuZyVOVUxMJTEsICAgICBSb25nJUNoZW4lTlVMTCUxLCAgICAgRmFuZ3lpbmclTHUlTlVMTCUxLCAgICAgWXVuZmVpJUx1JU5VTEwlMiwgICAgIFh1aHVpJUxpdSVOVUxMJTEsICAgICBZdXFpbmclQ2hlbiVOVUxMJTIsICAgICBZdXFpbmclQ2hlbiVOVUxMJTAsICAgICBYaW4lTGklTlVMTCUxLCAgICAgWW9uZyVMaSVOVUxMJTEsICAgICBIYW5zc2EgRHdhcmthJVN1bW1haCVOVUxMJTEsICAgICBIdWlodWFuZyVMaW4lTlVMTCUxLCAgICAgSmlheWFuZyVZYW4lTlVMTCUxLCAgICAgTWluJVpob3UlTlVMTCUwLCAgICAgSG9uZ3pob3UlTHUlTlVMTCUwLCAgICAgSG9uZ3pob3UlTHUlTlVMTCUwLCAgICAgSmllbWluZyVRdSVOVUxMJTAsICAgICBKaWVtaW5nJVF1JU5VTEwlMF0="
Set-B64Cell "I12" "X1BNQw=="
Set-B64Cell "E13" "W01pbmdsaSVZdWFuJU5VTEwlMSwgICAgIFdlbiVZaW4lTlVMTCUwLCAgICAgV2VuJVlpbiVOVUxMJTAsICAgICBaaGFvd3UlVGFvJU5VTEwlMSwgICAgIFdlaWp1biVUYW4lTlVMTCUxLCAgICAgWWklSHUlTlVMTCUwLCAgICAgT2xpdmVyJVNjaGlsZGdlbiVOVUxMJTIsICAgICBPbGl2ZXIlU2NoaWxkZ2VuJU5VTEwlMF0="
Set-B64Cell "I13" "X1BNQw=="
Set-B64Cell "E14" "W1Bpbmd6aGVuZyVNbyVOVUxMJTAsICAgICBZdWFueXVhbiVYaW5nJU5VTEwlMCwgICAgIFl1JVhpYW8lTlVMTCUwLCAgICAgTGlwaW5nJURlbmclTlVMTCUwLCAgICAgUWl1JVpoYW8lTlVMTCUwLCAgICAgSG9uZ2xpbmclV2FuZyVOVUxMJTAsICAgICBZb25nJVhpb25nJU5VTEwlMCwgICAgIFpoZW5zaHVuJUNoZW5nJU5VTEwlMCwgICAgIFNoaWNoZW5nJUdhbyVOVUxMJTAsICAgICBLZSVMaWFuZyVOVUxMJTAsICAgICBNaW5ncWklTHVvJU5VTEwlMCwgICAgIFRpZWxvbmclQ2hlbiVOVUxMJTAsICAgICBTaGlodWklU29uZyVOVUxMJTAsICAgICBaaGl5b25nJU1hJU5VTEwlMCwgICAgIFhpYW9waW5nJUNoZW4lTlVMTCUwLCAgICAgUnVpeWluZyVaaGVuZyVOVUxMJTAsICAgICBRaWFuJUNhbyVOVUxMJTAsICAgICBGYW4lV2FuZyVmYW5uZHl3YW5nQGZveG1haWwuY29tJTAsICAgICBZb25neGklWmhhbmclem5hY3QxOTM2QDEyNi5jb20lMF0="
Set-B64Cell "I14" "X1BNQw=="
Set-B64Cell "E15" "W0x1d2VuJVdhbmclTlVMTCUwLCAgICAgWHVuJUxpJU5VTEwlMSwgICAgIEh1aSVDaGVuJU5VTEwlMCwgICAgIFNoYW9uYW4lWWFuJU5VTEwlMSwgICAgIERvbmclTGklTlVMTCUxLCAgICAgWWFuJUxpJU5VTEwlMSwgICAgIFp1b2ppb25nJUdvbmclTlVMTCUxXQ=="
Set-B64Cell "I15" "X1BNQw=="
Set-B64Cell "E16" "W0d1cWluJVpoYW5nJU5VTEwlMCwgICAgIENoYW5nJUh1JU5VTEwlMSwgICAgIExpbmppZSVMdW8lTlVMTCUxLCAgICAgRmFuZyVGYW5nJU5VTEwlMSwgICAgIFlvbmdmZW5nJUNoZW4lTlVMTCUxLCAgICAgSmlhbmd1byVMaSVOVUxMJTEsICAgICBaaGl5b25nJVBlbmclTlVMTCUxLCAgICAgSHVhcWluJVBhbiVwaHEyMDEyQHdodS5lZHUuY24lMV0="
Set-B64Cell "I16" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E17" "W1RhbyVHdW8lTlVMTCUwLCAgICAgWW9uZ3poZW4lRmFuJU5VTEwlMSwgICAgIE1pbmclQ2hlbiVOVUxMJTIsICAgICBYaWFveWFuJVd1JU5VTEwlMiwgICAgIExpbiVaaGFuZyVOVUxMJTIsICAgICBUYW8lSGUlTlVMTCUyLCAgICAgSGFpcm9uZyVXYW5nJU5VTEwlMiwgICAgIEppbmclV2FuJU5VTEwlMSwgICAgIFhpbmdodWFuJVdhbmclTlVMTCUxLCAgICAgWmhpYmluZyVMdSVOVUxMJTJd"
Set-B64Cell "I17" "X1BNQw=="
Set-B64Cell "E18" "W0Jhcm5hYnklRC5QLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgQmVja2VyJUwuQi4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIENoZWxpY28lSi5ELiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgQ29oZW4lUy5MLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgQ29va2luZ2hhbSVKLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgQ29wcGElSy4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIENyYXdmb3JkJUouTS4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIERhdmlkc29uJUsuVy4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIERpZWZlbmJhY2glTS5BLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgRG9taW5lbGxvJUEuSi4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIER1ZXItSGVmZWxlJUouJWNvcmVHaXZlc05vRW1haWwlMCwgICBGYWx6b24lTC4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIEdpdGxpbiVKLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgSGFqaXphZGVoJU4uJWNvcmVHaXZlc05vRW1haWwlMCwgICBIYXJ2aW4lVC5HLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgSGlyc2NoJUouUy4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIEhpcnNjaHdlcmslRC5BLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgS2ltJUUuSi4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIEtvemVsJVouTS4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIE1hcnJhc3QlTC5NLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgTWNHaW5uJVQuJWNvcmVHaXZlc05vRW1haWwlMCwgICBNb2dhdmVybyVKLk4uJWNvcmVHaXZlc05vRW1haWwlMCwgICBOYXJhc2ltaGFuJU0uJWNvcmVHaXZlc05vRW1haWwlMCwgICBPc29yaW8lRy5BLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgUWl1JU0uJWNvcmVHaXZlc05vRW1haWwlMCwgICBSaWNoYXJkc29uJVMuJWNvcmVHaXZlc05vRW1haWwlMCwgICBaYW5vcyVULlAuJWNvcmVHaXZlc05vRW1haWwlMF0="
Set-B64Cell "F18" "bm90IGZvdW5k"
Set-B64Cell "G18" "Ti9B"
Set-B64Cell "I18" ""
Set-B64Cell "E19" "W05VTEwlTlVMTCVOVUxMJTAsICAgICBOVUxMJU5VTEwlTlVMTCUwLCAgICAgTlVMTCVOVUxMJU5VTEwlMCwgICAgIE5hbmN5JUNob3clTlVMTCUxLCAgICAgS2F0aGVyaW5lJUZsZW1pbmctRHV0cmElTlVMTCUxLCAgICAgUnlhbiVHaWVya2UlTlVMTCUxLCAgICAgQXJvbiVIYWxsJU5VTEwlMSwgICAgIE1pY2hlbGxlJUh1Z2hlcyVOVUxMJTEsICAgICBUYW1hcmElUGlsaXNodmlsaSVOVUxMJTEsICAgICBNYXR0aGV3JVJpdGNoZXklTlVMTCUxLCAgICAgS2F0aGVyaW5lJVJvZ3Vza2klTlVMTCUxLCAgICAgVGFtaSVTa29mZiVOVUxMJTEsICAgICBFbWlseSVVc3NlcnklTlVMTCUxXQ=="
Set-B64Cell "I19" "X1BNQw=="
Set-B64Cell "I20" ""
Set-B64Cell "F21" "bm90IGZvdW5k"
Set-B64Cell "G21" "Ti9B"
Set-B64Cell "I21" ""
Set-B64Cell "E22" "W0FudG9uZWxsaSVNYXNzaW1vJWNvcmVHaXZlc05vRW1haWwlMCwgICBDYWJyaW5pJUx1Y2ElY29yZUdpdmVzTm9FbWFpbCUwLCAgIENhc3RlbGxpJUFudG9uaW8lY29yZUdpdmVzTm9FbWFpbCUwLCAgIENlY2NvbmklTWF1cml6aW8lY29yZUdpdmVzTm9FbWFpbCUwLCAgIENlcmVkYSVEYW5pbG8lY29yZUdpdmVzTm9FbWFpbCUwLCAgIENvbHVjY2VsbG8lQW50b25pbyVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgRm90aSVHaXVzZXBwZSVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgRnVtYWdhbGxpJVJvYmVydG8lY29yZUdpdmVzTm9FbWFpbCUwLCAgIEdyYXNzZWxsaSVHaWFjb21vJWNvcmVHaXZlc05vRW1haWwlMCwgICBJb3R0aSVHaW9yZ2lvJWNvcmVHaXZlc05vRW1haWwlMCwgICBMYXRyb25pY28lTmljb2xhJWNvcmVHaXZlc05vRW1haWwlMCwgICBMb3JpbmklTHVjYSVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgTWVybGVyJVN0ZWZhbm8lY29yZUdpdmVzTm9FbWFpbCUwLCAgIE5hdGFsaW5pJUdpdXNlcHBlJWNvcmVHaXZlc05vRW1haWwlMCwgICBQZXNlbnRpJUFudG9uaW8lY29yZUdpdmVzTm9FbWFpbCUwLCAgIFBpYXR0aSVBbGVzc2FuZHJhJWNvcmVHaXZlc05vRW1haWwlMCwgICBSYW5pZXJpJU1hcmNvIFZpdG8lY29yZUdpdmVzTm9FbWFpbCUwLCAgIFNjYW5kcm9nbGlvJUFubmEgTWFyYSVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgU3RvcnRpJUVucmljbyVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgWmFuZWxsYSVBbGJlcnRvJWNvcmVHaXZlc05vRW1haWwlMCwgICBaYW5ncmlsbG8lQWxiZXJ0byVjb3JlR2l2ZXNOb0VtYWlsJTBd"
Set-B64Cell "F22" "bm90IGZvdW5k"
Set-B64Cell "G22" "Ti9B"
Set-B64Cell "I22" ""
Set-B64Cell "E23" "W1lpbmd6aGVuJUR1JU5VTEwlMCwgICAgIExlaSVUdSVOVUxMJTMsICAgICBMZWklVHUlTlVMTCUwLCAgICAgUGluZ2p1biVaaHUlTlVMTCUyLCAgICAgUGluZ2p1biVaaHUlTlVMTCUwLCAgICAgTWklTXUlTlVMTCUzLCAgICAgTWklTXUlTlVMTCUwLCAgICAgUnVuc2hlbmclV2FuZyVOVUxMJTIsICAgICBQZW5nY2hlbmclWWFuZyVOVUxMJTIsICAgICBYaSVXYW5nJU5VTEwlMSwgICAgIENoYW8lSHUlTlVMTCUyLCAgICAgUm9uZ3l1JVBpbmclTlVMTCUyLCAgICAgUGVuZyVIdSVOVUxMJTIsICAgICBUaWFuemhpJUxpJU5VTEwlMiwgICAgIEZlbmclQ2FvJU5VTEwlMSwgICAgIENocmlzdG9waGVyJUNoYW5nJU5VTEwlMSwgICAgIFFpbnlvbmclSHUlTlVMTCUyLCAgICAgWWFuZyVKaW4lTlVMTCUxLCAgICAgR3VvZ2FuZyVYdSVOVUxMJTJd"
Set-B64Cell "I23" "X1BNQw=="
Set-B64Cell "I24" "X1BNQw=="
Set-B64Cell "E25" "W05hbnNoYW4lQ2hlbiVOVUxMJTAsICAgICBNaW4lWmhvdSVOVUxMJTAsICAgICBYdWFuJURvbmclTlVMTCUwLCAgICAgSmllbWluZyVRdSVOVUxMJTAsICAgICBGZW5neXVuJUdvbmclTlVMTCUwLCAgICAgWWFuZyVIYW4lTlVMTCUwLCAgICAgWWFuZyVRaXUlTlVMTCUwLCAgICAgSmluZ2xpJVdhbmclTlVMTCUwLCAgICAgWWluZyVMaXUlTlVMTCUwLCAgICAgWXVhbiVXZWklTlVMTCUwLCAgICAgSmlhJ2FuJVhpYSVOVUxMJTAsICAgICBUaW5nJVl1JU5VTEwlMCwgICAgIFhpbnhpbiVaaGFuZyVOVUxMJTAsICAgICBMaSVaaGFuZyVOVUxMJTBd"
Set-B64Cell "I25" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E26" "W0NhcmJcdTAwZjMtQ2lzbmVybyVZYWNxdWVsaW4lY29yZUdpdmVzTm9FbWFpbCUwLCAgIEZlcm5cdTAwZTFuZGV6LUdvbnpcdTAwZTFsZXolUGF1bGElY29yZUdpdmVzTm9FbWFpbCUwLCAgIEhpZXJyZXp1ZWxvLVJvamFzJU5haWZpJWNvcmVHaXZlc05vRW1haWwlMCwgICBTdWJlcnQtU2FsYXMlTGl6YW5kcmElY29yZUdpdmVzTm9FbWFpbCUwXQ=="
Set-B64Cell "F26" "bm90IGZvdW5k"
Set-B64Cell "G26" "Ti9B"
Set-B64Cell "I26" ""
Set-B64Cell "I27" "X1BNQw=="
Set-B64Cell "E28" "W0NoYW9saW4lSHVhbmclTlVMTCUwLCAgICAgWWVtaW5nJVdhbmclTlVMTCUwLCAgICAgWGluZ3dhbmclTGklTlVMTCUwLCAgICAgTGlsaSVSZW4lTlVMTCUwLCAgICAgSmlhbnBpbmclWmhhbyVOVUxMJTAsICAgICBZaSVIdSVOVUxMJTAsICAgICBMaSVaaGFuZyVOVUxMJTAsICAgICBHdW9odWklRmFuJU5VTEwlMCwgICAgIEppdXlhbmclWHUlTlVMTCUwLCAgICAgWGlhb3lpbmclR3UlTlVMTCUwLCAgICAgWmhlbnNodW4lQ2hlbmclTlVMTCUwLCAgICAgVGluZyVZdSVOVUxMJTAsICAgICBKaWFhbiVYaWElTlVMTCUwLCAgICAgWXVhbiVXZWklTlVMTCUwLCAgICAgV2VuanVhbiVXdSVOVUxMJTAsICAgICBYdWVsZWklWGllJU5VTEwlMCwgICAgIFdlbiVZaW4lTlVMTCUwLCAgICAgSHVpJUxpJU5VTEwlMCwgICAgIE1pbiVMaXUlTlVMTCUwLCAgICAgWWFuJVhpYW8lTlVMTCUwLCAgICAgSG9uZyVHYW8lTlVMTCUwLCAgICAgTGklR3VvJU5VTEwlMCwgICAgIEp1bmdhbmclWGllJU5VTEwlMCwgICAgIEd1YW5nZmElV2FuZyVOVUxMJTAsICAgICBSb25nbWVuZyVKaWFuZyVOVUxMJTAsICAgICBaaGFuY2hlbmclR2FvJU5VTEwlMCwgICAgIFFpJUppbiVOVUxMJTAsICAgICBKaWFud2VpJVdhbmcld2FuZ2p3MjhAMTYzLmNvbSUwLCAgICAgQmluJUNhbyVjYW9iaW5fYmVuQDE2My5jb20lMF0="
Set-B64Cell "I28" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E29" "W0t1aSVMaXUlTlVMTCUwLCAgICAgWXVhbi1ZdWFuJUZhbmclTlVMTCUwLCAgICAgWWFuJURlbmclTlVMTCUwLCAgICAgV2VpJUxpdSVOVUxMJTAsICAgICBNZWktRmFuZyVXYW5nJU5VTEwlMCwgICAgIEppbmctUGluZyVNYSVOVUxMJTAsICAgICBXZWklWGlhbyVOVUxMJTAsICAgICBZaW5nLU5hbiVXYW5nJU5VTEwlMCwgICAgIE1pbi1IdWElWmhvbmclTlVMTCUwLCAgICAgQ2hlbmctSG9uZyVMaSVOVUxMJTAsICAgICBHdWFuZy1DYWklTGklTlVMTCUwLCAgICAgSHVpLUd1byVMaXUlTlVMTCUwLCAgICAgWGl1LVl1YW4lSGFvJU5VTEwlMCwgICAgIFhpdS1ZdWFuJUhhbyVOVUxMJTAsICAgICBQZWktRmFuZyVXZWklTlVMTCUwXQ=="
Set-B64Cell "I29" "X1BNQw=="
Set-B64Cell "E30" "W1lpbmd4aWElTGl1JXlpbmd4aWFsaXVAaG90bWFpbC5jb20lMCwgICAgIFlhbmclWWFuZyVOVUxMJTAsICAgICBDb25nJVpoYW5nJU5VTEwlMSwgICAgIEZlbmdtaW5nJUh1YW5nJU5VTEwlMSwgICAgIEZ1eGlhbmclV2FuZyVOVUxMJTEsICAgICBKaW5nJVl1YW4lTlVMTCUwLCAgICAgWmhhb3FpbiVXYW5nJU5VTEwlMSwgICAgIEppbnhpdSVMaSVOVUxMJTEsICAgICBKaWFubWluZyVMaSVOVUxMJTEsICAgICBDaGVuZyVGZW5nJU5VTEwlMSwgICAgIFpoZW5nJVpoYW5nJU5VTEwlMCwgICAgIExpZmVpJVdhbmclTlVMTCUxLCAgICAgTGluZyVQZW5nJU5VTEwlMSwgICAgIExpJUNoZW4lTlVMTCUxLCAgICAgWXVoYW8lUWluJU5VTEwlMSwgICAgIERhbmRhbiVaaGFvJU5VTEwlMSwgICAgIFNodWd1YW5nJVRhbiVOVUxMJTEsICAgICBMdSVZaW4lTlVMTCUxLCAgICAgSnVuJVh1JU5VTEwlMSwgICAgIENvbmd6aGFvJVpob3UlTlVMTCUxLCAgICAgQ2hlbmd5dSVKaWFuZyVqaWFuZ0BwdW1jLmVkdS5jbiUxLCAgICAgTGVpJUxpdSVsaXVsZWkzMzIyQGFsaXl1bi5jb20lMF0="
Set-B64Cell "I30" "X1BNQ19TcHJpbmdlcg=="
Set-B64Cell "E31" "W1dlaS1qaWUlR3VhbiVOVUxMJTAsICAgICBaaGVuZy15aSVOaSVOVUxMJTAsICAgICBaaGVuZy15aSVOaSVOVUxMJTAsICAgICBZdSVIdSVOVUxMJTAsICAgICBXZW4taHVhJUxpYW5nJU5VTEwlMCwgICAgIENodW4tcXVhbiVPdSVOVUxMJTAsICAgICBKaWFuLXhpbmclSGUlTlVMTCUwLCAgICAgTGVpJUxpdSVOVUxMJTAsICAgICBIb25nJVNoYW4lTlVMTCUwLCAgICAgQ2h1bi1saWFuZyVMZWklTlVMTCUwLCAgICAgRGF2aWQgUy5DLiVIdWklTlVMTCUwLCAgICAgQmluJUR1JU5VTEwlMCwgICAgIExhbi1qdWFuJUxpJU5VTEwlMCwgICAgIEd1YW5nJVplbmclTlVMTCUwLCAgICAgS3dvay1ZdW5nJVl1ZW4lTlVMTCUwLCAgICAgUnUtY2hvbmclQ2hlbiVOVUxMJTAsICAgICBDaHVuLWxpJVRhbmclTlVMTCUwLCAgICAgVGFvJVdhbmclTlVMTCUwLCAgICAgUGluZy15YW4lQ2hlbiVOVUxMJTAsICAgICBKaWUlWGlhbmclTlVMTCUwLCAgICAgU2hpLXl1ZSVMaSVOVUxMJTAsICAgICBKaW4tbGluJVdhbmclTlVMTCUwLCAgICAgWmktamluZyVMaWFuZyVOVUxMJTAsICAgICBZaS14aWFuZyVQZW5nJU5VTEwlMCwgICAgIExpJVdlaSVOVUxMJTAsICAgICBZb25nJUxpdSVOVUxMJTAsICAgICBZYS1odWElSHUlTlVMTCUwLCAgICAgUGVuZyVQZW5nJU5VTEwlMCwgICAgIEppYW4tbWluZyVXYW5nJU5VTEwlMCwgICAgIEppLXlhbmclTGl1JU5VTEwlMCwgICAgIFpob25nJUNoZW4lTlVMTCUwLCAgICAgR2FuZyVMaSVOVUxMJTAsICAgICBaaGktamlhbiVaaGVuZyVOVUxMJTAsICAgICBTaGFvLXFpbiVRaXUlTlVMTCUwLCAgICAgSmllJUx1byVOVUxMJTAsICAgICBDaGFuZy1qaWFuZyVZZSVOVUxMJTAsICAgICBTaGFvLXlvbmclWmh1JU5VTEwlMCwgICAgIE5hbi1zaGFuJVpob25nJU5VTEwlMF0="
Set-B64Cell "I31" "X1BNQw=="
Set-B64Cell "E32" "W1JvbmctSHVpJUR1JU5VTEwlMCwgICAgIExpLU1pbiVMaXUlTlVMTCUwLCAgICAgV2VuJVlpbiVOVUxMJTAsICAgICBXZW4lV2FuZyVOVUxMJTAsICAgICBMdS1MdSVHdWFuJU5VTEwlMCwgICAgIE1pbmctTGklWXVhbiVOVUxMJTAsICAgICBZdS1MZWklTGklTlVMTCUwLCAgICAgWWklSHUlTlVMTCUwLCAgICAgWHUtWWFuJUxpJU5VTEwlMCwgICAgIEJpbmclU3VuJU5VTEwlMCwgICAgIFBlbmclUGVuZyVOVUxMJTAsICAgICBIdWFuLVpob25nJVNoaSVOVUxMJTBd"
Set-B64Cell "I32" "X1BNQw=="
Set-B64Cell "E33" "W1dlaS1qaWUlR3VhbiVOVUxMJTAsICAgICBXZW4taHVhJUxpYW5nJU5VTEwlMCwgICAgIFlpJVpoYW8lTlVMTCUyLCAgICAgSGVuZy1ydWklTGlhbmclTlVMTCUyLCAgICAgWmktc2hlbmclQ2hlbiVOVUxMJTIsICAgICBZaS1taW4lTGklTlVMTCUyLCAgICAgWGlhby1xaW5nJUxpdSVOVUxMJTIsICAgICBSdS1jaG9uZyVDaGVuJU5VTEwlMCwgICAgIENodW4tbGklVGFuZyVOVUxMJTAsICAgICBUYW8lV2FuZyVOVUxMJTAsICAgICBDaHVuLXF1YW4lT3UlTlVMTCUwLCAgICAgTGklTGklTlVMTCUwLCAgICAgUGluZy15YW4lQ2hlbiVOVUxMJTAsICAgICBMaW5nJVNhbmclTlVMTCU0LCAgICAgV2VpJVdhbmclTlVMTCUwLCAgICAgSmlhbi1mdSVMaSVOVUxMJTIsICAgICBDYWktY2hlbiVMaSVOVUxMJTIsICAgICBMaS1taW4lT3UlTlVMTCUyLCAgICAgQm8lQ2hlbmclTlVMTCUyLCAgICAgU2hhbiVYaW9uZyVOVUxMJTIsICAgICBaaGVuZy15aSVOaSVOVUxMJTAsICAgICBKaWUlWGlhbmclTlVMTCUwLCAgICAgWXUlSHUlTlVMTCUwLCAgICAgTGVpJUxpdSVOVUxMJTAsICAgICBIb25nJVNoYW4lTlVMTCUwLCAgICAgQ2h1bi1saWFuZyVMZWklTlVMTCUwLCAgICAgWWkteGlhbmclUGVuZyVOVUxMJTAsICAgICBMaSVXZWklTlVMTCUwLCAgICAgWW9uZyVMaXUlTlVMTCUwLCAgICAgWWEtaHVhJUh1JU5VTEwlMCwgICAgIFBlbmclUGVuZyVOVUxMJTAsICAgICBKaWFuLW1pbmclV2FuZyVOVUxMJTAsICAgICBKaS15YW5nJUxpdSVOVUxMJTAsICAgICBaaG9uZyVDaGVuJU5VTEwlMCwgICAgIEdhbmclTGklTlVMTCUwLCAgICAgWmhpLWppYW4lWmhlbmclTlVMTCUwLCAgICAgU2hhby1xaW4lUWl1JU5VTEwlMCwgICAgIEppZSVMdW8lTlVMTCUwLCAgICAgQ2hhbmctamlhbmclWWUlTlVMTCUwLCAgICAgU2hhby15b25nJVpodSVOVUxMJTAsICAgICBMaW4tbGluZyVDaGVuZyVOVUxMJTIsICAgICBGZW5nJVllJU5VTEwlNCwgICAgIFNoaS15dWUlTGklTlVMTCUwLCAgICAgSmluLXBpbmclWmhlbmclTlVMTCUyLCAgICAgTnVvLWZ1JVpoYW5nJU5VTEwlMiwgICAgIE5hbi1zaGFuJVpob25nJU5VTEwlMCwgICAgIEppYW4teGluZyVIZSVOVUxMJTBd"
Set-B64Cell "I33" "X1BNQw=="
Set-B64Cell "E34" "W1hpYW8tV2VpJVh1JU5VTEwlMCwgICAgIFhpYW8tWGluJVd1JU5VTEwlMCwgICAgIFhpYW4tR2FvJUppYW5nJU5VTEwlMCwgICAgIEthaS1KaW4lWHUlTlVMTCUwLCAgICAgTGluZy1KdW4lWWluZyVOVUxMJTAsICAgICBDaHVuLUxpYW4lTWElTlVMTCUwLCAgICAgU2hpLUJvJUxpJU5VTEwlMCwgICAgIEh1YS1ZaW5nJVdhbmclTlVMTCUwLCAgICAgU2hlbmclWmhhbmclTlVMTCUwLCAgICAgSGFpLU52JUdhbyVOVUxMJTAsICAgICBKaS1GYW5nJVNoZW5nJU5VTEwlMCwgICAgIEhvbmctTGl1JUNhaSVOVUxMJTAsICAgICBZdW4tUWluZyVRaXUlTlVMTCUwLCAgICAgTGFuLUp1YW4lTGklTlVMTCUwXQ=="
Set-B64Cell "I34" "X1BNQw=="
Set-B64Cell "E35" "W1RpZUxvbmclQ2hlbiVOVUxMJTAsICAgICBaaGUlRGFpJU5VTEwlMSwgICAgIFBpbmd6aGVuZyVNbyVOVUxMJTEsICAgICBYaW55dSVMaSVOVUxMJTEsICAgICBaaGl5b25nJU1hJU5VTEwlMCwgICAgIFNoaWh1aSVTb25nJU5VTEwlMCwgICAgIFhpYW9waW5nJUNoZW4lTlVMTCUwLCAgICAgTWluZ3FpJUx1byVOVUxMJTAsICAgICBLZSVMaWFuZyVOVUxMJTAsICAgICBTaGljaGVuZyVHYW8lTlVMTCUwLCAgICAgWW9uZ3hpJVpoYW5nJU5VTEwlMCwgICAgIExpcGluZyVEZW5nJWRlbmdkZW5nNzhAMTI2LmNvbSUwLCAgICAgWW9uZyVYaW9uZyVOVUxMJTAsICAgICBZb25nJVhpb25nJU5VTEwlMF0="
Set-B64Cell "I35" "X1BNQw=="
Set-B64Cell "E36" "W1dlbi1odWElTGlhbmclTlVMTCUwLCAgICAgV2VpLWppZSVHdWFuJU5VTEwlMSwgICAgIENhaS1jaGVuJUxpJU5VTEwlMCwgICAgIFlpLW1pbiVMaSVOVUxMJTAsICAgICBIZW5nLXJ1aSVMaWFuZyVOVUxMJTAsICAgICBZaSVaaGFvJU5VTEwlMCwgICAgIFhpYW8tcWluZyVMaXUlTlVMTCUwLCAgICAgTGluZyVTYW5nJU5VTEwlMCwgICAgIFJ1LWNob25nJUNoZW4lTlVMTCUwLCAgICAgQ2h1bi1saSVUYW5nJU5VTEwlMCwgICAgIFRhbyVXYW5nJU5VTEwlMCwgICAgIFdlaSVXYW5nJU5VTEwlMCwgICAgIFFpLWh1YSVIZSVOVUxMJTEsICAgICBaaS1zaGVuZyVDaGVuJU5VTEwlMCwgICAgIFNvb2stU2FuJVdvbmclTlVMTCUxLCAgICAgTWFyayVaYW5pbiVOVUxMJTEsICAgICBKdW4lTGl1JU5VTEwlMCwgICAgIFhpbiVYdSVOVUxMJTAsICAgICBKdW4lSHVhbmclTlVMTCUxLCAgICAgSmlhbi1mdSVMaSVOVUxMJTAsICAgICBMaS1taW4lT3UlTlVMTCUwLCAgICAgQm8lQ2hlbmclTlVMTCUwLCAgICAgU2hhbiVYaW9uZyVOVUxMJTAsICAgICBaaGFuLWhvbmclWGllJU5VTEwlMSwgICAgIFpoZW5nLXlpJU5pJU5VTEwlMCwgICAgIFl1JUh1JU5VTEwlMCwgICAgIExlaSVMaXUlTlVMTCUwLCAgICAgSG9uZyVTaGFuJU5VTEwlMCwgICAgIENodW4tbGlhbmclTGVpJU5VTEwlMCwgICAgIFlpLXhpYW5nJVBlbmclTlVMTCUwLCAgICAgTGklV2VpJU5VTEwlMCwgICAgIFlvbmclTGl1JU5VTEwlMCwgICAgIFlhLWh1YSVIdSVOVUxMJTAsICAgICBQZW5nJVBlbmclTlVMTCUwLCAgICAgSmlhbi1taW5nJVdhbmclTlVMTCUwLCAgICAgSmkteWFuZyVMaXUlTlVMTCUwLCAgICAgWmhvbmclQ2hlbiVOVUxMJTAsICAgICBHYW5nJUxpJU5VTEwlMCwgICAgIFpoaS1qaWFuJVpoZW5nJU5VTEwlMCwgICAgIFNoYW8tcWluJVFpdSVOVUxMJTAsICAgICBKaWUlTHVvJU5VTEwlMCwgICAgIENoYW5nLWppYW5nJVllJU5VTEwlMCwgICAgIFNoYW8teW9uZyVaaHUlTlVMTCUwLCAgICAgTGluLWxpbmclQ2hlbmclTlVMTCUwLCAgICAgRmVuZyVZZSVOVUxMJTAsICAgICBTaGkteXVlJUxpJU5VTEwlMCwgICAgIEppbi1waW5nJVpoZW5nJU5VTEwlMCwgICAgIE51by1mdSVaaGFuZyVOVUxMJTAsICAgICBOYW4tc2hhbiVaaG9uZyVOVUxMJTAsICAgICBKaWFuLXhpbmclSGUlTlVMTCUwXQ=="
Set-B64Cell "I36" "X1BNQw=="
Set-B64Cell "E37" "W0xhbmclV2FuZyVOVUxMJTAsICAgICBXZW5ibyVIZSVOVUxMJTIsICAgICBYaWFvbWVpJVl1JU5VTEwlMiwgICAgIERhbG9uZyVIdSVOVUxMJTIsICAgICBNaW5nd2VpJUJhbyVOVUxMJTIsICAgICBIdWFmZW4lTGl1JU5VTEwlMiwgICAgIEppYWxpJVpob3UlTlVMTCUyLCAgICAgSG9uZyVKaWFuZyVOVUxMJTJd"
Set-B64Cell "I37" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E38" "W0ppYW5sZWklQ2FvJU5VTEwlMCwgICAgIFdlbi1KdW4lVHUldHV3ZW5qdW5AaXJtLWNhbXMuYWMuY24lMCwgICAgIFdlbmxpbiVDaGVuZyVOVUxMJTAsICAgICBMZWklWXUlTlVMTCUwLCAgICAgWWEtS3VuJUxpdSVOVUxMJTAsICAgICBYaWFveW9uZyVIdSVOVUxMJTAsICAgICBRaWFuZyVMaXUlTlVMTCUwXQ=="
Set-B64Cell "I38" "X1BNQw=="
Set-B64Cell "E39" "W0ZlaSVaaG91JU5VTEwlMCwgICAgIFRpbmclWXUlTlVMTCUwLCAgICAgUm9uZ2h1aSVEdSVOVUxMJTAsICAgICBHdW9odWklRmFuJU5VTEwlMCwgICAgIFlpbmclTGl1JU5VTEwlMCwgICAgIFpoaWJvJUxpdSVOVUxMJTAsICAgICBKaWUlWGlhbmclTlVMTCUwLCAgICAgWWVtaW5nJVdhbmclTlVMTCUwLCAgICAgQmluJVNvbmclTlVMTCUwLCAgICAgWGlhb3lpbmclR3UlTlVMTCUwLCAgICAgTHVsdSVHdWFuJU5VTEwlMCwgICAgIFl1YW4lV2VpJU5VTEwlMCwgICAgIEh1aSVMaSVOVUxMJTAsICAgICBYdWRvbmclV3UlTlVMTCUwLCAgICAgSml1eWFuZyVYdSVOVUxMJTAsICAgICBTaGVuZ2ppbiVUdSVOVUxMJTAsICAgICBZaSVaaGFuZyVOVUxMJTAsICAgICBIdWElQ2hlbiVOVUxMJTAsICAgICBCaW4lQ2FvJU5VTEwlMF0="
Set-B64Cell "I39" "X1BNQ19lbHNldmllcg=="
Set-B64Cell "E40" "W1J1aSVaaGFuZyVOVUxMJTAsICAgICBIdWFuZ3FpbmclT3V5YW5nJU5VTEwlMSwgICAgIExpbmdsaSVGdSVOVUxMJTEsICAgICBTaGlqaWUlV2FuZyVOVUxMJTEsICAgICBKaWFuZ2xvbmclSGFuJU5VTEwlMSwgICAgIEtlamllJUh1YW5nJU5VTEwlMSwgICAgIE1pbmdmYW5nJUppYSVOVUxMJTEsICAgICBRaWJpbiVTb25nJU5VTEwlMSwgICAgIFpoZW5taW5nJUZ1JWRhdmlkZnV6bWluZ0B3aHUuZWR1LmNuJTFd"
Set-B64Cell "I40" "X1BNQ19TcHJpbmdlcg=="
Set-B64Cell "E41" "W1RhbyVDaGVuJU5VTEwlMCwgICAgIERpJVd1JU5VTEwlMywgICAgIEh1aWxvbmclQ2hlbiVOVUxMJTMsICAgICBXZWltaW5nJVlhbiVOVUxMJTMsICAgICBEYW5sZWklWWFuZyVOVUxMJTMsICAgICBHdWFuZyVDaGVuJU5VTEwlMywgICAgIEtlJU1hJU5VTEwlMywgICAgIERvbmclWHUlTlVMTCU1LCAgICAgSGFpamluZyVZdSVOVUxMJTMsICAgICBIb25nd3UlV2FuZyVOVUxMJTMsICAgICBUYW8lV2FuZyVOVUxMJTAsICAgICBXZWklR3VvJU5VTEwlMywgICAgIEppYSVDaGVuJU5VTEwlMywgICAgIENoZW4lRGluZyVOVUxMJTMsICAgICBYaWFvcGluZyVaaGFuZyVOVUxMJTMsICAgICBKaWFxdWFuJUh1YW5nJU5VTEwlMywgICAgIE1laWZhbmclSGFuJU5VTEwlMywgICAgIFNodXNoZW5nJUxpJU5VTEwlMCwgICAgIFhpYW9waW5nJUx1byVOVUxMJTMsICAgICBKaWFucGluZyVaaGFvJU5VTEwlMCwgICAgIFFpbiVOaW5nJU5VTEwlM10="
Set-B64Cell "I41" "X1BNQw=="
